# Automatische test-sync: 2025-06-19 15:30:10
$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append new row 16 with the unsubscribe message ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A16").Value = "Afmelding nieuwsbrief"
$ws.Range("B16").Value = "mailmind.test@zohomail.eu"
$ws.Range("C16").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$ws.Range("D16").Value = "Afmelding"
$ws.Range("F16").Value = "2025-06-19 15:28:11"
$ws.Range("G16").Value = "Nee"

# Extend the conditional formatting ranges so they cover the new row too.
$dRules = $ws.Range("D2:D15").FormatConditions
$dRules.Item(1).ModifyAppliesToRange($ws.Range("D2:D16"))

$gRules = $ws.Range("G2:G15").FormatConditions
$gRules.Item(1).ModifyAppliesToRange($ws.Range("G2:G16"))

# --- "Dashboard" sheet: bump the "Afmelding" counter from 2 to 3 ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 3
